{"js": "// The Good-Faith-Estimate notice paragraph currently reads:\n//   \"You have a right to receive the Good Faith Estimate in writing at least\n//    one day before scheduling an appointment for any non-emergency medical\n//    services.\"\n// The edit removes the phrase \" in writing\" so it reads:\n//   \"You have a right to receive the Good Faith Estimate at least one day\n//    before scheduling an appointment for any non-emergency medical\n//    services.\"\n\nconst target = \"Good Faith Estimate in writing\";\nconst replacement = \"Good Faith Estimate\";\n\nconst results = context.document.body.search(target, { matchCase: true, matchWholeWord: false });\nresults.load(\"items,text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(`Could not find the text \"${target}\" to update.`);\n}\n\n// Replace every match (expected to be exactly one) in place, preserving the\n// surrounding formatting of the run(s) it lives in.\nfor (const r of results.items) {\n  r.insertText(replacement, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The Good-Faith-Estimate notice paragraph currently reads:\n#   \"You have a right to receive the Good Faith Estimate in writing at least\n#    one day before scheduling an appointment for any non-emergency medical\n#    services.\"\n# The edit removes the phrase \" in writing\" so it reads:\n#   \"You have a right to receive the Good Faith Estimate at least one day\n#    before scheduling an appointment for any non-emergency medical\n#    services.\"\n\n$d = $word.ActiveDocument\n\n$find = \"Good Faith Estimate in writing\"\n$replace = \"Good Faith Estimate\"\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Replacement.ClearFormatting()\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$rng.Find.Execute($find, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replace, $wdReplaceAll)\n"}
